# Generate Report for Handoff
# The "b.md" file has moved from "Handed back: in sync with en-US" to
# "Ready for handoff" with a fresh handoff file/timestamp. Update the
# Overview sheet plus the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the "b.md" row.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-22 02:29:16"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the "b.md" row.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-22 02:29:12"

foreach ($h in $zhcn.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the "b.md" row.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-22 02:29:16"

foreach ($h in $dede.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
